$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Addr, $Val)
    $r = $ws.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "34.203.86"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.789.72"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.16%  "
Set-TextValue "D5" "226.33"
$ws.Range("E5").Value = "  -0.25%  "
Set-TextValue "D6" "0.549"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  +0.13%  "
Set-TextValue "D8" "32.40"
$ws.Range("E8").Value = "  +0.60%  "
Set-TextValue "D9" "0.294"
$ws.Range("E9").Value = "  +0.17%  "
Set-TextValue "D10" "0.0690"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "2.046.92"
$ws.Range("E12").Value = "  +0.05%  "
Set-TextValue "D13" "11.15"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "1.785.21"
$ws.Range("E14").Value = "  -0.44%  "
Set-TextValue "D15" "0.627"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "34.182.71"
Set-TextValue "D17" "4.20"
$ws.Range("E17").Value = "  +0.55%  "
Set-TextValue "D18" "67.93"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0807"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "246.40"
$ws.Range("E20").Value = "  +1.06%  "
Set-TextValue "D21" "11.05"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  +0.24%  "
Set-TextValue "D23" "4.18"
$ws.Range("E23").Value = "  +2.06%  "
Set-TextValue "D24" "2.05"
$ws.Range("E24").Value = "  +0.50%  "
Set-TextValue "D25" "162.00"
$ws.Range("E25").Value = "  +0.07%  "
Set-TextValue "D26" "7.18"
$ws.Range("E26").Value = "  -0.27%  "
Set-TextValue "D27" "16.32"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.0522"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "1.23"
$ws.Range("E31").Value = "  -0.18%  "
Set-TextValue "D32" "3.78"
$ws.Range("E32").Value = "  +3.64%  "
Set-TextValue "D33" "3.74"
$ws.Range("E33").Value = "  +3.65%  "
Set-TextValue "D34" "1.81"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "1.445.38"
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("E36").Value = "  +9.66%  "
Set-TextValue "D37" "0.663"
$ws.Range("E37").Value = "  +2.60%  "
Set-TextValue "D38" "0.0191"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  +1.19%  "
Set-TextValue "D40" "82.27"
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "0.924"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "13.65"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("E45").Value = "  +1.76%  "
Set-TextValue "D46" "6.13"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "1.945.80"
$ws.Range("E48").Value = "  -0.04%  "
Set-TextValue "D49" "105.30"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  -6.61%  "
